$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.620274999999999
$ws.Range("H2").Value = 22.860825
$ws.Range("I2").Value = 0.6584612850834004
$ws.Range("J2").Value = 0.6584612850834003
$ws.Range("M2").Value = 36.81180933333333
$ws.Range("N2").Value = 110.435428
$ws.Range("O2").Value = 0.2598784967371026
$ws.Range("P2").Value = 0.2598784967371026
$ws.Range("Q2").Value = 280.5161103675666
$ws.Range("R2").Value = 2524.644993308099
$ws.Range("S2").Value = 0.1711199289270549
$ws.Range("T2").Value = 0.1711199289270548
$ws.Range("G3").Value = 7.620274999999999
$ws.Range("H3").Value = 22.860825
$ws.Range("I3").Value = 0.6584612850834004
$ws.Range("J3").Value = 0.6584612850834003
$ws.Range("O3").Value = 0.1970278712683331
$ws.Range("P3").Value = 0.197027871268333
$ws.Range("Q3").Value = 212.6743565786666
$ws.Range("R3").Value = 1914.069209208
$ws.Range("S3").Value = 0.1297352253125934
$ws.Range("T3").Value = 0.1297352253125933
$ws.Range("G4").Value = 7.620274999999999
$ws.Range("H4").Value = 22.860825
$ws.Range("I4").Value = 0.6584612850834004
$ws.Range("J4").Value = 0.6584612850834003
$ws.Range("M4").Value = 21.95609833333333
$ws.Range("N4").Value = 65.868295
$ws.Range("O4").Value = 0.1550023737603119
$ws.Range("P4").Value = 0.1550023737603119
$ws.Range("Q4").Value = 167.3115072270417
$ws.Range("R4").Value = 1505.803565043375
$ws.Range("S4").Value = 0.1020630622171925
$ws.Range("T4").Value = 0.1020630622171925
$ws.Range("G5").Value = 7.620274999999999
$ws.Range("H5").Value = 22.860825
$ws.Range("I5").Value = 0.6584612850834004
$ws.Range("J5").Value = 0.6584612850834003
$ws.Range("M5").Value = 13.23098133333333
$ws.Range("N5").Value = 39.692944
$ws.Range("O5").Value = 0.09340609987756826
$ws.Range("P5").Value = 0.09340609987756825
$ws.Range("Q5").Value = 100.8237162798667
$ws.Range("R5").Value = 907.4134465187999
$ws.Range("S5").Value = 0.06150430056001205
$ws.Range("T5").Value = 0.06150430056001203
$ws.Range("G6").Value = 7.620274999999999
$ws.Range("H6").Value = 22.860825
$ws.Range("I6").Value = 0.6584612850834004
$ws.Range("J6").Value = 0.6584612850834003
$ws.Range("M6").Value = 22.080681
$ws.Range("N6").Value = 66.242043
$ws.Range("O6").Value = 0.1558818838066577
$ws.Range("P6").Value = 0.1558818838066577
$ws.Range("Q6").Value = 168.260861407275
$ws.Range("R6").Value = 1514.347752665475
$ws.Range("S6").Value = 0.1026421855325531
$ws.Range("T6").Value = 0.1026421855325531
$ws.Range("G7").Value = 7.620274999999999
$ws.Range("H7").Value = 22.860825
$ws.Range("I7").Value = 0.6584612850834004
$ws.Range("J7").Value = 0.6584612850834003
$ws.Range("M7").Value = 19.66149466666667
$ws.Range("N7").Value = 58.984484
$ws.Range("O7").Value = 0.1388032745500265
$ws.Range("P7").Value = 0.1388032745500265
$ws.Range("Q7").Value = 149.8259962710333
$ws.Range("R7").Value = 1348.4339664393
$ws.Range("S7").Value = 0.09139658253399449
$ws.Range("T7").Value = 0.09139658253399448
$ws.Range("I8").Value = 0.262323813236933
$ws.Range("J8").Value = 0.262323813236933
$ws.Range("M8").Value = 36.81180933333333
$ws.Range("N8").Value = 110.435428
$ws.Range("O8").Value = 0.2598784967371026
$ws.Range("P8").Value = 0.2598784967371026
$ws.Range("Q8").Value = 111.7545669168569
$ws.Range("R8").Value = 1005.791102251712
$ws.Range("S8").Value = 0.06817231824235861
$ws.Range("T8").Value = 0.06817231824235861
$ws.Range("I9").Value = 0.262323813236933
$ws.Range("J9").Value = 0.262323813236933
$ws.Range("O9").Value = 0.1970278712683331
$ws.Range("P9").Value = 0.197027871268333
$ws.Range("S9").Value = 0.05168510250506469
$ws.Range("T9").Value = 0.05168510250506468
$ws.Range("I10").Value = 0.262323813236933
$ws.Range("J10").Value = 0.262323813236933
$ws.Range("M10").Value = 21.95609833333333
$ws.Range("N10").Value = 65.868295
$ws.Range("O10").Value = 0.1550023737603119
$ws.Range("P10").Value = 0.1550023737603119
$ws.Range("Q10").Value = 66.65508446507556
$ws.Range("R10").Value = 599.89576018568
$ws.Range("S10").Value = 0.04066081374558136
$ws.Range("T10").Value = 0.04066081374558136
$ws.Range("I11").Value = 0.262323813236933
$ws.Range("J11").Value = 0.262323813236933
$ws.Range("M11").Value = 13.23098133333333
$ws.Range("N11").Value = 39.692944
$ws.Range("O11").Value = 0.09340609987756826
$ws.Range("P11").Value = 0.09340609987756825
$ws.Range("Q11").Value = 40.16707180575288
$ws.Range("R11").Value = 361.503646251776
$ws.Range("S11").Value = 0.02450264429947353
$ws.Range("T11").Value = 0.02450264429947353
$ws.Range("I12").Value = 0.262323813236933
$ws.Range("J12").Value = 0.262323813236933
$ws.Range("M12").Value = 22.080681
$ws.Range("N12").Value = 66.242043
$ws.Range("O12").Value = 0.1558818838066577
$ws.Range("P12").Value = 0.1558818838066577
$ws.Range("Q12").Value = 67.033296843408
$ws.Range("R12").Value = 603.299671590672
$ws.Range("S12").Value = 0.04089153017471898
$ws.Range("T12").Value = 0.04089153017471898
$ws.Range("I13").Value = 0.262323813236933
$ws.Range("J13").Value = 0.262323813236933
$ws.Range("M13").Value = 19.66149466666667
$ws.Range("N13").Value = 58.984484
$ws.Range("O13").Value = 0.1388032745500265
$ws.Range("P13").Value = 0.1388032745500265
$ws.Range("Q13").Value = 59.68904710754844
$ws.Range("R13").Value = 537.2014239679361
$ws.Range("S13").Value = 0.03641140426973589
$ws.Range("T13").Value = 0.03641140426973589
$ws.Range("G14").Value = 0.9167423333333334
$ws.Range("H14").Value = 2.750227
$ws.Range("I14").Value = 0.07921490167966665
$ws.Range("J14").Value = 0.07921490167966663
$ws.Range("M14").Value = 36.81180933333333
$ws.Range("N14").Value = 110.435428
$ws.Range("O14").Value = 0.2598784967371026
$ws.Range("P14").Value = 0.2598784967371026
$ws.Range("Q14").Value = 33.74694398246178
$ws.Range("R14").Value = 303.722495842156
$ws.Range("S14").Value = 0.02058624956768915
$ws.Range("T14").Value = 0.02058624956768915
$ws.Range("G15").Value = 0.9167423333333334
$ws.Range("H15").Value = 2.750227
$ws.Range("I15").Value = 0.07921490167966665
$ws.Range("J15").Value = 0.07921490167966663
$ws.Range("O15").Value = 0.1970278712683331
$ws.Range("P15").Value = 0.197027871268333
$ws.Range("Q15").Value = 25.58537400423111
$ws.Range("R15").Value = 230.26836603808
$ws.Range("S15").Value = 0.01560754345067502
$ws.Range("T15").Value = 0.01560754345067502
$ws.Range("G16").Value = 0.9167423333333334
$ws.Range("H16").Value = 2.750227
$ws.Range("I16").Value = 0.07921490167966665
$ws.Range("J16").Value = 0.07921490167966663
$ws.Range("M16").Value = 21.95609833333333
$ws.Range("N16").Value = 65.868295
$ws.Range("O16").Value = 0.1550023737603119
$ws.Range("P16").Value = 0.1550023737603119
$ws.Range("Q16").Value = 20.12808481699611
$ws.Range("R16").Value = 181.152763352965
$ws.Range("S16").Value = 0.01227849779753805
$ws.Range("T16").Value = 0.01227849779753805
$ws.Range("G17").Value = 0.9167423333333334
$ws.Range("H17").Value = 2.750227
$ws.Range("I17").Value = 0.07921490167966665
$ws.Range("J17").Value = 0.07921490167966663
$ws.Range("M17").Value = 13.23098133333333
$ws.Range("N17").Value = 39.692944
$ws.Range("O17").Value = 0.09340609987756826
$ws.Range("P17").Value = 0.09340609987756825
$ws.Range("Q17").Value = 12.12940069980978
$ws.Range("R17").Value = 109.164606298288
$ws.Range("S17").Value = 0.007399155018082693
$ws.Range("T17").Value = 0.00739915501808269
$ws.Range("G18").Value = 0.9167423333333334
$ws.Range("H18").Value = 2.750227
$ws.Range("I18").Value = 0.07921490167966665
$ws.Range("J18").Value = 0.07921490167966663
$ws.Range("M18").Value = 22.080681
$ws.Range("N18").Value = 66.242043
$ws.Range("O18").Value = 0.1558818838066577
$ws.Range("P18").Value = 0.1558818838066577
$ws.Range("Q18").Value = 20.242295021529
$ws.Range("R18").Value = 182.180655193761
$ws.Range("S18").Value = 0.01234816809938561
$ws.Range("T18").Value = 0.01234816809938561
$ws.Range("G19").Value = 0.9167423333333334
$ws.Range("H19").Value = 2.750227
$ws.Range("I19").Value = 0.07921490167966665
$ws.Range("J19").Value = 0.07921490167966663
$ws.Range("M19").Value = 19.66149466666667
$ws.Range("N19").Value = 58.984484
$ws.Range("O19").Value = 0.1388032745500265
$ws.Range("P19").Value = 0.1388032745500265
$ws.Range("Q19").Value = 149.8259962710333
$ws.Range("R19").Value = 1348.4339664393
$ws.Range("S19").Value = 0.09139658253399449
$ws.Range("T19").Value = 0.09139658253399448
